# Regenerate orders with updated distance/sizes.
# The edit replaces several text tokens that appear throughout the
# "Condition", "Filename_Left", "Filename_Right", "Distance" and "Size"
# columns (and the corresponding shared-string table entries once saved):
#   D80 -> D86
#   D64 -> D69
#   D51 -> D55
#   S30 -> S31
# None of the replacement tokens re-introduce one of the other source
# tokens, so the four replacements can be performed independently, in any
# order, across every used cell in the sheet.

$wb = $excel.ActiveWorkbook

$xlWhole = 2      # xlPart = 1, xlWhole = 2
$xlByRows = 1      # xlByColumns = 2, xlByRows = 1
$falseVal = $false

foreach ($ws in $wb.Worksheets) {
    $rng = $ws.UsedRange

    $rng.Replace("D80", "D86", $xlWhole, $xlByRows, $falseVal, $falseVal, $falseVal) | Out-Null
    $rng.Replace("D64", "D69", $xlWhole, $xlByRows, $falseVal, $falseVal, $falseVal) | Out-Null
    $rng.Replace("D51", "D55", $xlWhole, $xlByRows, $falseVal, $falseVal, $falseVal) | Out-Null
    $rng.Replace("S30", "S31", $xlWhole, $xlByRows, $falseVal, $falseVal, $falseVal) | Out-Null
}
